# Apply the scraped-symbol-list refresh (Wed Jan  4 05:45:01 UTC 2023 run).
# Price/volume/name/link cells are stored as *text* (t="inlineStr") in the
# workbook, not as numbers - so every write below goes through Set-TextCell,
# which prefixes the value with a quote-prefix apostrophe (forces Excel to
# keep it as literal text instead of silently parsing "255.39" or "3.95%"
# into a number/percentage) and then resets the cell Style back to "Normal"
# so no stray number-format/quote-prefix style sticks around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Set-PlainCell($addr, $text) {
    $ws.Range($addr).Value = $text
}

# Row 2
Set-TextCell "D2" '255.39'
Set-TextCell "E2" '3.95%'

# Row 3
Set-TextCell "D3" '28.02'
Set-TextCell "E3" '-3.84%'

# Row 4
Set-TextCell "D4" '5.359'
Set-TextCell "E4" '4.22%'

# Row 5
Set-TextCell "D5" '0.05822'
Set-TextCell "E5" '0.62%'

# Row 6
Set-TextCell "D6" '6.713'
Set-TextCell "E6" '1.49%'

# Row 7
Set-PlainCell "B7" 'GateToken'
Set-PlainCell "C7" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell "D7" '3.225'
Set-TextCell "E7" '1.56%'

# Row 8
Set-PlainCell "B8" 'MXToken'
Set-PlainCell "C8" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell "D8" '0.8649'
Set-TextCell "E8" '0.84%'

# Row 9
Set-PlainCell "B9" 'FTXToken'
Set-PlainCell "C9" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell "D9" '0.9135'
Set-TextCell "E9" '5.92%'

# Row 10
Set-PlainCell "B10" 'One'
Set-PlainCell "C10" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell "D10" '0.01061'
Set-TextCell "E10" '1,671.10%'

# Row 11
Set-PlainCell "B11" 'WazirX'
Set-PlainCell "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell "D11" '0.1424'
Set-TextCell "E11" '4.10%'

# Row 12
Set-PlainCell "B12" 'MandalaExchangeToken'
Set-PlainCell "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell "D12" '0.07178'
Set-TextCell "E12" '1.69%'

# Row 13
Set-PlainCell "B13" 'BitrueCoin'
Set-PlainCell "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell "D13" '0.03203'
Set-TextCell "E13" '-2.66%'

# Row 14
Set-PlainCell "B14" 'BitMartToken'
Set-PlainCell "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell "D14" '0.09239'
Set-TextCell "E14" '-1.28%'

# Row 15
Set-PlainCell "B15" 'BitForexToken'
Set-PlainCell "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell "D15" '0.001548'
Set-TextCell "E15" '1.70%'

# Row 16
Set-PlainCell "B16" 'TigerCash'
Set-PlainCell "C16" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell "D16" '0.005904'
Set-TextCell "E16" '-1.19%'

# Row 17
Set-PlainCell "B17" 'LEO'
Set-PlainCell "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell "D17" '3.496'
Set-TextCell "E17" '0.30%'

# Row 18
Set-TextCell "E18" '4.05%'

# Row 19
Set-TextCell "E19" '-1.01%'

# Row 20
Set-TextCell "D20" '0.03443'
Set-TextCell "E20" '3.35%'

# Row 21
Set-TextCell "D21" '0.1334'
Set-TextCell "E21" '3.92%'

# Row 22
Set-TextCell "D22" '3.534'
Set-TextCell "E22" '11.07%'

# Row 23
Set-TextCell "D23" '0.04157'
Set-TextCell "E23" '0.53%'

# Row 24
Set-TextCell "E24" '-1.38%'

# Row 25
Set-PlainCell "B25" 'HotbitToken'
Set-PlainCell "C25" 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell "D25" '0.005127'
Set-TextCell "E25" '23.77%'

# Row 26
Set-PlainCell "B26" 'BitKan'
Set-PlainCell "C26" 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell "D26" '0.001223'
Set-TextCell "E26" '-0.12%'

# Row 27
Set-TextCell "E27" '-0.73%'

# Row 28
Set-TextCell "D28" '0.0001938'
Set-TextCell "E28" '34.24%'

# Row 40
Set-TextCell "D40" '0.03844'
Set-TextCell "E40" '3.04%'

# Row 41
Set-TextCell "D41" '0.005737'
Set-TextCell "E41" '-0.28%'

# Row 42
Set-TextCell "D42" '0.1099'
Set-TextCell "E42" '2.77%'

# Row 43
Set-TextCell "D43" '0.002200'
Set-TextCell "E43" '0.10%'

# Row 44
Set-TextCell "D44" '0.009891'
Set-TextCell "E44" '7.79%'

# Row 45
Set-TextCell "D45" '0.00005288'
Set-TextCell "E45" '0.09%'

# Row 46
Set-TextCell "E46" '0.10%'

# Row 47
Set-TextCell "E47" '72.76%'

# Row 48
Set-TextCell "D48" '0.002211'
Set-TextCell "E48" '1.77%'

# Row 49
Set-TextCell "D49" '0.00002100'
Set-TextCell "E49" '0.10%'

# Row 50
Set-TextCell "D50" '0.0002000'
Set-TextCell "E50" '0.10%'
